# Generate Report for Handoff
# Adds two new files (5413f0c7-d550-4aac-8f4a-9a7bd51e688b, 6fe4f66e-300b-4d1f-93a3-6cd180100228)
# to the localization-status workbook, each "Ready for handoff", ahead of the
# existing 96e8afa5-9aa7-4dfe-8212-60b1e5cc62e8 entry (which shifts down).

$wb = $excel.ActiveWorkbook

$HYPER_COLOR = 15570276   # OLE BGR int for RGB FF6495ED (the sheet's HyperLink font color)
$DATETIME_FMT = "yyyy-mm-dd HH:mm:ss"

function Style-AsLink($rng) {
    $rng.Font.Underline = $true
    $rng.Font.Color = $HYPER_COLOR
}

function Style-AsDateTime($rng) {
    $rng.NumberFormat = $DATETIME_FMT
}

# ---------------------------------------------------------------------------
# Sheet "Overview": columns File Name | zh-cn | de-de | Latest Handoff Date
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

# Row 3: now the 5413f0c7 entry (was 96e8afa5)
$wsOverview.Range("A3").Value = "5413f0c7-d550-4aac-8f4a-9a7bd51e688b.md"
$wsOverview.Hyperlinks.Add($wsOverview.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/ec9f525bab77160e4b2b9b74fa9582a7c7521877/e2e/5413f0c7-d550-4aac-8f4a-9a7bd51e688b.md", "", "", "5413f0c7-d550-4aac-8f4a-9a7bd51e688b.md") | Out-Null
Style-AsLink $wsOverview.Range("A3")
$wsOverview.Range("B3").Value = "Ready for handoff"
$wsOverview.Range("C3").Value = "Ready for handoff"
$wsOverview.Range("D3").Value = "2016-03-22 00:35:10"
Style-AsDateTime $wsOverview.Range("D3")

# Row 4: 6fe4f66e entry (new)
$wsOverview.Range("A4").Value = "6fe4f66e-300b-4d1f-93a3-6cd180100228.md"
$wsOverview.Hyperlinks.Add($wsOverview.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/4b82b1586d2bf781c8fa90c99a512c454efa324d/e2e/6fe4f66e-300b-4d1f-93a3-6cd180100228.md", "", "", "6fe4f66e-300b-4d1f-93a3-6cd180100228.md") | Out-Null
Style-AsLink $wsOverview.Range("A4")
$wsOverview.Range("B4").Value = "Ready for handoff"
$wsOverview.Range("C4").Value = "Ready for handoff"
$wsOverview.Range("D4").Value = "2016-03-22 00:35:10"
Style-AsDateTime $wsOverview.Range("D4")

# Row 5: 96e8afa5 entry (was row 3, pushed down)
$wsOverview.Range("A5").Value = "96e8afa5-9aa7-4dfe-8212-60b1e5cc62e8.md"
$wsOverview.Hyperlinks.Add($wsOverview.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/3073ab0b57de045a493efbd6ee1c3f4aaaaaa857/e2e/96e8afa5-9aa7-4dfe-8212-60b1e5cc62e8.md", "", "", "96e8afa5-9aa7-4dfe-8212-60b1e5cc62e8.md") | Out-Null
Style-AsLink $wsOverview.Range("A5")
$wsOverview.Range("B5").Value = "Ready for handoff"
$wsOverview.Range("C5").Value = "Ready for handoff"
$wsOverview.Range("D5").Value = "2016-03-22 00:33:37"
Style-AsDateTime $wsOverview.Range("D5")

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

# Row 3: 5413f0c7 entry
$wsZh.Range("A3").Value = "5413f0c7-d550-4aac-8f4a-9a7bd51e688b.md"
$wsZh.Hyperlinks.Add($wsZh.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/ec9f525bab77160e4b2b9b74fa9582a7c7521877/e2e/5413f0c7-d550-4aac-8f4a-9a7bd51e688b.md", "", "", "5413f0c7-d550-4aac-8f4a-9a7bd51e688b.md") | Out-Null
Style-AsLink $wsZh.Range("A3")
$wsZh.Range("B3").Value = ".md"
$wsZh.Range("C3").Value = "Ready for handoff"
$wsZh.Range("D3").Value = "5413f0c7-d550-4aac-8f4a-9a7bd51e688b.ec9f525bab77160e4b2b9b74fa9582a7c7521877.zh-cn.xlf"
$wsZh.Hyperlinks.Add($wsZh.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ec9f525bab77160e4b2b9b74fa9582a7c7521877/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/5413f0c7-d550-4aac-8f4a-9a7bd51e688b.ec9f525bab77160e4b2b9b74fa9582a7c7521877.zh-cn.xlf", "", "", "5413f0c7-d550-4aac-8f4a-9a7bd51e688b.ec9f525bab77160e4b2b9b74fa9582a7c7521877.zh-cn.xlf") | Out-Null
Style-AsLink $wsZh.Range("D3")
$wsZh.Range("E3").Value = "2016-03-22 00:35:06"
Style-AsDateTime $wsZh.Range("E3")
$wsZh.Range("H3").Value = "0001-01-01 00:00:00"
Style-AsDateTime $wsZh.Range("H3")
$wsZh.Range("J3").Value = "Include"

# Row 4: 6fe4f66e entry
$wsZh.Range("A4").Value = "6fe4f66e-300b-4d1f-93a3-6cd180100228.md"
$wsZh.Hyperlinks.Add($wsZh.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/4b82b1586d2bf781c8fa90c99a512c454efa324d/e2e/6fe4f66e-300b-4d1f-93a3-6cd180100228.md", "", "", "6fe4f66e-300b-4d1f-93a3-6cd180100228.md") | Out-Null
Style-AsLink $wsZh.Range("A4")
$wsZh.Range("B4").Value = ".md"
$wsZh.Range("C4").Value = "Ready for handoff"
$wsZh.Range("D4").Value = "6fe4f66e-300b-4d1f-93a3-6cd180100228.4b82b1586d2bf781c8fa90c99a512c454efa324d.zh-cn.xlf"
$wsZh.Hyperlinks.Add($wsZh.Range("D4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/4b82b1586d2bf781c8fa90c99a512c454efa324d/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/6fe4f66e-300b-4d1f-93a3-6cd180100228.4b82b1586d2bf781c8fa90c99a512c454efa324d.zh-cn.xlf", "", "", "6fe4f66e-300b-4d1f-93a3-6cd180100228.4b82b1586d2bf781c8fa90c99a512c454efa324d.zh-cn.xlf") | Out-Null
Style-AsLink $wsZh.Range("D4")
$wsZh.Range("E4").Value = "2016-03-22 00:35:06"
Style-AsDateTime $wsZh.Range("E4")
$wsZh.Range("H4").Value = "0001-01-01 00:00:00"
Style-AsDateTime $wsZh.Range("H4")
$wsZh.Range("J4").Value = "Include"

# Row 5: 96e8afa5 entry (was row 3, pushed down)
$wsZh.Range("A5").Value = "96e8afa5-9aa7-4dfe-8212-60b1e5cc62e8.md"
$wsZh.Hyperlinks.Add($wsZh.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/3073ab0b57de045a493efbd6ee1c3f4aaaaaa857/e2e/96e8afa5-9aa7-4dfe-8212-60b1e5cc62e8.md", "", "", "96e8afa5-9aa7-4dfe-8212-60b1e5cc62e8.md") | Out-Null
Style-AsLink $wsZh.Range("A5")
$wsZh.Range("B5").Value = ".md"
$wsZh.Range("C5").Value = "Ready for handoff"
$wsZh.Range("D5").Value = "96e8afa5-9aa7-4dfe-8212-60b1e5cc62e8.0012e40d796e5c6f54b3c87d5af7bf616b8ae37b.zh-cn.xlf"
$wsZh.Hyperlinks.Add($wsZh.Range("D5"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/a8d307b6cb8b29c0798d457611105d1993f1f720/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/96e8afa5-9aa7-4dfe-8212-60b1e5cc62e8.0012e40d796e5c6f54b3c87d5af7bf616b8ae37b.zh-cn.xlf", "", "", "96e8afa5-9aa7-4dfe-8212-60b1e5cc62e8.0012e40d796e5c6f54b3c87d5af7bf616b8ae37b.zh-cn.xlf") | Out-Null
Style-AsLink $wsZh.Range("D5")
$wsZh.Range("E5").Value = "2016-03-22 00:33:33"
Style-AsDateTime $wsZh.Range("E5")
$wsZh.Range("H5").Value = "0001-01-01 00:00:00"
Style-AsDateTime $wsZh.Range("H5")
$wsZh.Range("J5").Value = "Include"

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

# Row 3: 5413f0c7 entry
$wsDe.Range("A3").Value = "5413f0c7-d550-4aac-8f4a-9a7bd51e688b.md"
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/ec9f525bab77160e4b2b9b74fa9582a7c7521877/e2e/5413f0c7-d550-4aac-8f4a-9a7bd51e688b.md", "", "", "5413f0c7-d550-4aac-8f4a-9a7bd51e688b.md") | Out-Null
Style-AsLink $wsDe.Range("A3")
$wsDe.Range("B3").Value = ".md"
$wsDe.Range("C3").Value = "Ready for handoff"
$wsDe.Range("D3").Value = "5413f0c7-d550-4aac-8f4a-9a7bd51e688b.ec9f525bab77160e4b2b9b74fa9582a7c7521877.de-de.xlf"
$wsDe.Hyperlinks.Add($wsDe.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ec9f525bab77160e4b2b9b74fa9582a7c7521877/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/5413f0c7-d550-4aac-8f4a-9a7bd51e688b.ec9f525bab77160e4b2b9b74fa9582a7c7521877.de-de.xlf", "", "", "5413f0c7-d550-4aac-8f4a-9a7bd51e688b.ec9f525bab77160e4b2b9b74fa9582a7c7521877.de-de.xlf") | Out-Null
Style-AsLink $wsDe.Range("D3")
$wsDe.Range("E3").Value = "2016-03-22 00:35:10"
Style-AsDateTime $wsDe.Range("E3")
$wsDe.Range("H3").Value = "0001-01-01 00:00:00"
Style-AsDateTime $wsDe.Range("H3")
$wsDe.Range("J3").Value = "Include"

# Row 4: 6fe4f66e entry
$wsDe.Range("A4").Value = "6fe4f66e-300b-4d1f-93a3-6cd180100228.md"
$wsDe.Hyperlinks.Add($wsDe.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/4b82b1586d2bf781c8fa90c99a512c454efa324d/e2e/6fe4f66e-300b-4d1f-93a3-6cd180100228.md", "", "", "6fe4f66e-300b-4d1f-93a3-6cd180100228.md") | Out-Null
Style-AsLink $wsDe.Range("A4")
$wsDe.Range("B4").Value = ".md"
$wsDe.Range("C4").Value = "Ready for handoff"
$wsDe.Range("D4").Value = "6fe4f66e-300b-4d1f-93a3-6cd180100228.4b82b1586d2bf781c8fa90c99a512c454efa324d.de-de.xlf"
$wsDe.Hyperlinks.Add($wsDe.Range("D4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/4b82b1586d2bf781c8fa90c99a512c454efa324d/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/6fe4f66e-300b-4d1f-93a3-6cd180100228.4b82b1586d2bf781c8fa90c99a512c454efa324d.de-de.xlf", "", "", "6fe4f66e-300b-4d1f-93a3-6cd180100228.4b82b1586d2bf781c8fa90c99a512c454efa324d.de-de.xlf") | Out-Null
Style-AsLink $wsDe.Range("D4")
$wsDe.Range("E4").Value = "2016-03-22 00:35:10"
Style-AsDateTime $wsDe.Range("E4")
$wsDe.Range("H4").Value = "0001-01-01 00:00:00"
Style-AsDateTime $wsDe.Range("H4")
$wsDe.Range("J4").Value = "Include"

# Row 5: 96e8afa5 entry (was row 3, pushed down)
$wsDe.Range("A5").Value = "96e8afa5-9aa7-4dfe-8212-60b1e5cc62e8.md"
$wsDe.Hyperlinks.Add($wsDe.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/3073ab0b57de045a493efbd6ee1c3f4aaaaaa857/e2e/96e8afa5-9aa7-4dfe-8212-60b1e5cc62e8.md", "", "", "96e8afa5-9aa7-4dfe-8212-60b1e5cc62e8.md") | Out-Null
Style-AsLink $wsDe.Range("A5")
$wsDe.Range("B5").Value = ".md"
$wsDe.Range("C5").Value = "Ready for handoff"
$wsDe.Range("D5").Value = "96e8afa5-9aa7-4dfe-8212-60b1e5cc62e8.0012e40d796e5c6f54b3c87d5af7bf616b8ae37b.de-de.xlf"
$wsDe.Hyperlinks.Add($wsDe.Range("D5"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/335346405c8c72e7301abad97d5697d2151ff791/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/96e8afa5-9aa7-4dfe-8212-60b1e5cc62e8.0012e40d796e5c6f54b3c87d5af7bf616b8ae37b.de-de.xlf", "", "", "96e8afa5-9aa7-4dfe-8212-60b1e5cc62e8.0012e40d796e5c6f54b3c87d5af7bf616b8ae37b.de-de.xlf") | Out-Null
Style-AsLink $wsDe.Range("D5")
$wsDe.Range("E5").Value = "2016-03-22 00:33:37"
Style-AsDateTime $wsDe.Range("E5")
$wsDe.Range("H5").Value = "0001-01-01 00:00:00"
Style-AsDateTime $wsDe.Range("H5")
$wsDe.Range("J5").Value = "Include"

Write-Host "Handoff report rows added."
